# Applies the crypto price/volume/ranking update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new values (plain decimal-looking numbers) would be
# auto-converted to numeric cells by Excel's normal type inference,
# losing precision and changing the cell type from text to number.
# The source column is formatted as text, so force these specific
# cells to Text format first, preserving them as exact strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "91.850.21"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "3.121.04"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "242.71"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "625.45"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E7").Value = "  +6.02%  "
$ws.Range("D8").Value = "0.375"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.121.66"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "0.775"
$ws.Range("E11").Value = "  +6.81%  "
$ws.Range("E12").Value = "  +3.33%  "
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "35.62"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").Value = "91.709.52"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "5.51"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "3.706.58"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "3.131.37"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").Value = "14.77"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "5.84"
$ws.Range("E22").Value = "  +3.09%  "
$ws.Range("D23").Value = "448.39"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "5.90"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "91.87"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "11.94"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("D28").Value = "3.297.79"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "0.252"
$ws.Range("E30").Value = "  +24.27%  "
$ws.Range("E31").Value = "  +14.68%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "9.26"
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.116"
$ws.Range("E33").Value = "  +34.36%  "
$ws.Range("E34").Value = "  +23.52%  "
$ws.Range("E35").Value = "  +10.76%  "
$ws.Range("D36").Value = "26.66"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("D37").Value = "7.61"
$ws.Range("E37").Value = "  +5.52%  "
$ws.Range("D38").Value = "4.13"
$ws.Range("E38").Value = "  +20.91%  "
$ws.Range("E39").Value = "  -6.01%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "492.89"
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "22.17"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "156.69"
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("D48").Value = "0.694"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").Value = "4.57"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "44.76"
$ws.Range("E51").Value = "  -2.55%  "
